$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 6000
$ws.Range("I76").Value = 6000
$ws.Range("K76").Value = 6000
$ws.Range("M76").Value = -5685

$ws.Range("H79").Value = 6000
$ws.Range("I79").Value = 6000
$ws.Range("K79").Value = 6000
$ws.Range("M79").Value = -4908

$ws.Range("H98").Value = 1843.9166
$ws.Range("I98").Value = 1553.1
$ws.Range("J98").Value = 3298
$ws.Range("K98").Value = 1553.1
$ws.Range("L98").Value = 3298
$ws.Range("M98").Value = -55.09999999999991
$ws.Range("N98").Value = -6294

$ws.Range("H99").Value = 387.25
$ws.Range("I99").Value = 387.25
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1161.75
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 336.25
$ws.Range("N99").ClearContents()

$ws.Range("H101").Value = 883.3333
$ws.Range("I101").Value = 325
$ws.Range("K101").Value = 975
$ws.Range("M101").Value = 647

$ws.Range("H122").Value = 1843.9166
$ws.Range("I122").Value = 1553.1
$ws.Range("J122").Value = 3298
$ws.Range("K122").Value = 4659.299999999999
$ws.Range("L122").Value = 9894
$ws.Range("M122").Value = -2209.299999999999
$ws.Range("N122").Value = -14794

$ws.Range("H125").Value = 2100
$ws.Range("I125").Value = 2765.5
$ws.Range("J125").Value = 1656.3334
$ws.Range("K125").Value = 24889.5
$ws.Range("L125").Value = 14907.0006
$ws.Range("M125").Value = -22429.5
$ws.Range("N125").Value = -19827.0006

$ws.Range("H132").Value = 1305.0741
$ws.Range("I132").Value = 1214.4783
$ws.Range("J132").Value = 1826
$ws.Range("K132").Value = 3643.4349
$ws.Range("L132").Value = 5478
$ws.Range("M132").Value = -1113.4349
$ws.Range("N132").Value = -10538

$ws.Range("H135").Value = 13023
$ws.Range("I135").Value = 14998
$ws.Range("J135").Value = 12035.5
$ws.Range("K135").Value = 134982
$ws.Range("L135").Value = 108319.5
$ws.Range("M135").Value = -132447
$ws.Range("N135").Value = -113389.5

$ws.Range("H137").Value = 8367.575999999999
$ws.Range("I137").Value = 4440.5483
$ws.Range("K137").Value = 13321.6449
$ws.Range("M137").Value = -10771.6449

$ws.Range("H138").Value = 3487.1528
$ws.Range("I138").Value = 2248.7632
$ws.Range("J138").Value = 4871.2354
$ws.Range("K138").Value = 6746.2896
$ws.Range("L138").Value = 14613.7062
$ws.Range("M138").Value = -1606.2896
$ws.Range("N138").Value = -24893.7062

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4973.125
$ws.Range("I2").Value = 4071.4119
$ws.Range("K2").Value = 4071.4119
$ws.Range("M2").Value = -3958.4119

$ws.Range("H45").Value = 11120.643
$ws.Range("I45").Value = 12557.417
$ws.Range("K45").Value = 12557.417
$ws.Range("M45").Value = -12180.417

$ws.Range("H74").Value = 8232.762000000001
$ws.Range("I74").Value = 7682.0938
$ws.Range("K74").Value = 7682.0938
$ws.Range("M74").Value = -6808.0938

$ws.Range("H77").Value = 8232.762000000001
$ws.Range("I77").Value = 7682.0938
$ws.Range("K77").Value = 38410.469
$ws.Range("M77").Value = -34042.469

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H116").Value = 4973.125
$ws.Range("I116").Value = 4071.4119
$ws.Range("K116").Value = 4071.4119
$ws.Range("M116").Value = -1777.4119

$ws.Range("H130").Value = 61405
$ws.Range("J130").Value = 61405
$ws.Range("L130").Value = 61405
$ws.Range("N130").Value = -71445

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4973.125
$ws.Range("I3").Value = 4071.4119
$ws.Range("K3").Value = 4071.4119
$ws.Range("M3").Value = -3957.4119

$ws.Range("H58").Value = 44999.332
$ws.Range("I58").Value = 39999
$ws.Range("K58").Value = 39999
$ws.Range("M58").Value = -39705

$ws.Range("H94").Value = 1597.9697
$ws.Range("I94").Value = 1540.129
$ws.Range("J94").Value = 2494.5
$ws.Range("K94").Value = 1540.129
$ws.Range("L94").Value = 2494.5
$ws.Range("M94").Value = -1089.129
$ws.Range("N94").Value = -3396.5

$ws.Range("H134").Value = 5849.172
$ws.Range("I134").Value = 4454.597
$ws.Range("K134").Value = 13363.791
$ws.Range("M134").Value = -10828.791

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 34988.285
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 34988.285
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 34988.285
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -35328.285

$ws.Range("H31").Value = 1503.4546
$ws.Range("I31").Value = 839.5
$ws.Range("K31").Value = 839.5
$ws.Range("M31").Value = -544.5

$ws.Range("H34").Value = 1503.4546
$ws.Range("I34").Value = 839.5
$ws.Range("K34").Value = 839.5
$ws.Range("M34").Value = -637.5

$ws.Range("H58").Value = 2637.7021
$ws.Range("J58").Value = 5331.7144
$ws.Range("L58").Value = 5331.7144
$ws.Range("N58").Value = -5737.7144

$ws.Range("H96").Value = 24539.666
$ws.Range("J96").Value = 24539.666
$ws.Range("L96").Value = 24539.666
$ws.Range("N96").Value = -30031.666

$ws.Range("H116").Value = 65042
$ws.Range("J116").Value = 65042
$ws.Range("L116").Value = 65042
$ws.Range("N116").Value = -74220

$ws.Range("H136").Value = 2637.7021
$ws.Range("J136").Value = 5331.7144
$ws.Range("L136").Value = 15995.1432
$ws.Range("N136").Value = -21095.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1919.6666
$ws.Range("I5").Value = 878
$ws.Range("J5").Value = 2029.3158
$ws.Range("K5").Value = 2634
$ws.Range("L5").Value = 6087.9474
$ws.Range("M5").Value = -2522
$ws.Range("N5").Value = -6311.9474

$ws.Range("H68").Value = 1920.625
$ws.Range("I68").Value = 1346.25
$ws.Range("K68").Value = 4038.75
$ws.Range("M68").Value = -3227.75

$ws.Range("H71").Value = 1920.625
$ws.Range("I71").Value = 1346.25
$ws.Range("K71").Value = 12116.25
$ws.Range("M71").Value = -8060.25

$ws.Range("H121").Value = 1774.0938
$ws.Range("J121").Value = 2039.0834
$ws.Range("L121").Value = 6117.2502
$ws.Range("N121").Value = -8737.2502

$ws.Range("H135").Value = 1919.6666
$ws.Range("I135").Value = 878
$ws.Range("J135").Value = 2029.3158
$ws.Range("K135").Value = 7902
$ws.Range("L135").Value = 18263.8422
$ws.Range("M135").Value = -5367
$ws.Range("N135").Value = -23333.8422

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

$ws.Range("H140").Value = 8846.454
$ws.Range("I140").Value = 8924.111000000001
$ws.Range("J140").Value = 8497
$ws.Range("K140").Value = 26772.333
$ws.Range("L140").Value = 25491
$ws.Range("M140").Value = -21592.333
$ws.Range("N140").Value = -35851

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H42").Value = 89999
$ws.Range("J42").Value = 89999
$ws.Range("L42").Value = 89999
$ws.Range("N42").Value = -90969

$ws.Range("H113").Value = 252131.56
$ws.Range("I113").Value = 323098
$ws.Range("K113").Value = 323098
$ws.Range("M113").Value = -320928

$ws.Range("H115").Value = 89999
$ws.Range("J115").Value = 89999
$ws.Range("L115").Value = 89999
$ws.Range("N115").Value = -92349

$ws.Range("H126").Value = 4087.6667
$ws.Range("I126").Value = 4084.1428
$ws.Range("K126").Value = 12252.4284
$ws.Range("M126").Value = -9782.428400000001

$ws.Range("H130").Value = 90778
$ws.Range("J130").Value = 90778
$ws.Range("L130").Value = 90778
$ws.Range("N130").Value = -100818

$ws.Range("H132").Value = 8810.919
$ws.Range("I132").Value = 9978.414000000001
$ws.Range("J132").Value = 4578.75
$ws.Range("K132").Value = 29935.242
$ws.Range("L132").Value = 13736.25
$ws.Range("M132").Value = -27405.242
$ws.Range("N132").Value = -18796.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1460.1936
$ws.Range("I46").Value = 1066.3334
$ws.Range("J46").Value = 1554.72
$ws.Range("K46").Value = 1066.3334
$ws.Range("L46").Value = 1554.72
$ws.Range("M46").Value = -878.3334
$ws.Range("N46").Value = -1930.72

$ws.Range("H61").Value = 2998.7273
$ws.Range("I61").Value = 2998.7273
$ws.Range("K61").Value = 2998.7273
$ws.Range("M61").Value = -2796.7273

$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H113").Value = 2998.7273
$ws.Range("I113").Value = 2998.7273
$ws.Range("K113").Value = 2998.7273
$ws.Range("M113").Value = -828.7273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H115").Value = 399999
$ws.Range("J115").Value = 399999
$ws.Range("L115").Value = 399999
$ws.Range("N115").Value = -403133

$ws.Range("H132").Value = 21011.303
$ws.Range("I132").Value = 12980.086
$ws.Range("J132").Value = 34396.668
$ws.Range("K132").Value = 38940.258
$ws.Range("L132").Value = 103190.004
$ws.Range("M132").Value = -36410.258
$ws.Range("N132").Value = -108250.004
